$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdaten")

# New single-column content (A1=Suchbegriff, A2=Testmanager, A3=Testanalyst)
$ws.Range("A1").Value = "Suchbegriff"
$ws.Range("A2").Value = "Testmanager"
$ws.Range("A3").Value = "Testanalyst"

# Drop the highlight/fill style that used to sit on A2:A3
$ws.Range("A1:A3").ClearFormats()

# Remove the old second column (URL / FreelancerMap / FreelanceDE)
$ws.Range("B1:B3").ClearContents()
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws.Columns.Item(2).Delete()

$ws.Range("A1:A1048576").Select()
